# Applies the edits described by the commit:
#  - Slide 2: "Frontend - Vercel" -> "Frontend - Vercel (html / css / js)"
#  - Slide 2: "Fly.io / Cloud Run / Railway" -> "Cloud Run (FastAPI)"
#  - Slide 4: merge "Reference: ... Algorithm - " + "GeeksforGeeks" runs into one run
#  - Slide 4: fix double space in "is further  described" -> "is further described"

$enDash = [char]0x2013

$p = $ppt.ActivePresentation

# --- Slide 2: Frontend / Backend labels (nested inside "Group 12") ---
$s2 = $p.Slides.Item(2)
$grp12 = $s2.Shapes.Item(10)
$archBox = $grp12.GroupItems.Item(2)
$archTr = $archBox.TextFrame.TextRange

$archFull = $archTr.Text
$frontendIdx = $archFull.IndexOf("Frontend")
$backendIdx = $archFull.IndexOf("Backend")
$frontendLen = $backendIdx - $frontendIdx
$frontendChars = $archTr.Characters($frontendIdx + 1, $frontendLen)
$frontendChars.Text = "Frontend " + $enDash + " Vercel (html / css / js)"

$archFull2 = $archTr.Text
$backendValIdx = $archFull2.IndexOf("Fly.io / Cloud Run / Railway")
$backendValChars = $archTr.Characters($backendValIdx + 1, "Fly.io / Cloud Run / Railway".Length)
$backendValChars.Text = "Cloud Run (FastAPI)"

# --- Slide 4: BM25 reference + GeeksforGeeks merged into a single run ---
$s4 = $p.Slides.Item(4)
$bm25Shape = $s4.Shapes.Item(4)
$bm25Tr = $bm25Shape.TextFrame.TextRange

$bm25Full = $bm25Tr.Text
$refIdx = $bm25Full.IndexOf("Reference:")
$refLen = $bm25Full.Length - $refIdx
$refChars = $bm25Tr.Characters($refIdx + 1, $refLen)
$refChars.Text = "Reference:  What is BM25 (Best Matching 25) Algorithm " + $enDash + " GeeksforGeeks"

# --- Slide 4: remove double space in "is further  described" ---
$downstreamShape = $s4.Shapes.Item(5)
$downstreamTr = $downstreamShape.TextFrame.TextRange

$downstreamFull = $downstreamTr.Text
$furtherIdx = $downstreamFull.IndexOf("Further downstream")
$furtherLen = $downstreamFull.Length - $furtherIdx
$furtherChars = $downstreamTr.Characters($furtherIdx + 1, $furtherLen)
$furtherChars.Text = "Further downstream scoring using the 2 components (dense + sparse) is further described in separate slide"
